# Tastyworks Trading workbook update - 2-Jun-2021 midday update
$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Name = "Trades"
$wb.Worksheets.Item("Sheet2").Name = "Account U27637"

$wsTrades  = $wb.Worksheets.Item("Trades")
$wsAccount = $wb.Worksheets.Item("Account U27637")

Write-Output "done"
